$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: each entry maps row number to column letter/value pairs that changed
$rowUpdates = @(
    @{ Row = 2; D = '28.179.08'; E = '  +0.64%  ' },
    @{ Row = 3; D = '1.878.78'; E = '  +1.09%  ' },
    @{ Row = 4; D = '1.007'; E = '  +0.29%  ' },
    @{ Row = 5; D = '314.29'; E = '  +0.86%  ' },
    @{ Row = 6; D = '1.006'; E = '  +0.28%  ' },
    @{ Row = 7; D = '0.5136'; E = '  +0.95%  ' },
    @{ Row = 8; D = '0.3906'; E = '  +2.52%  ' },
    @{ Row = 9; E = '  +1.45%  ' },
    @{ Row = 10; E = '  +1.12%  ' },
    @{ Row = 11; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '6.229'; E = '  +0.63%  ' },
    @{ Row = 12; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '20.66'; E = '  +1.20%  ' },
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.873.11'; E = '  +0.46%  ' },
    @{ Row = 14; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '7.263'; E = '  +1.23%  ' },
    @{ Row = 15; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.005'; E = '  +0.03%  ' },
    @{ Row = 16; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.00001101'; E = '  +0.73%  ' },
    @{ Row = 17; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '91.27'; E = '  +0.96%  ' },
    @{ Row = 18; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.06679'; E = '  +1.04%  ' },
    @{ Row = 19; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '17.80'; E = '  +0.48%  ' },
    @{ Row = 20; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.005'; E = '  +0.25%  ' },
    @{ Row = 21; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '6.045'; E = '  +0.65%  ' },
    @{ Row = 22; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '28.207.72'; E = '  +0.69%  ' },
    @{ Row = 23; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '11.13'; E = '  +1.01%  ' },
    @{ Row = 24; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '2.264'; E = '  +1.61%  ' },
    @{ Row = 25; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.089.94'; E = '  +0.58%  ' },
    @{ Row = 26; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '159.80'; E = '  +1.36%  ' },
    @{ Row = 27; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.482'; E = '  -2.78%  ' },
    @{ Row = 28; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '20.74'; E = '  +1.86%  ' },
    @{ Row = 29; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '125.23'; E = '  +0.62%  ' },
    @{ Row = 30; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1062'; E = '  +0.33%  ' },
    @{ Row = 31; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '1.039'; E = '  +0.44%  ' },
    @{ Row = 32; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '5.848'; E = '  +4.58%  ' },
    @{ Row = 33; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '3.611'; E = '  +0.24%  ' },
    @{ Row = 34; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '9.631'; E = '  +0.15%  ' },
    @{ Row = 35; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.02447'; E = '  +1.67%  ' },
    @{ Row = 36; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.06564'; E = '  +0.60%  ' },
    @{ Row = 37; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.2189'; E = '  +0.90%  ' },
    @{ Row = 38; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.199'; E = '  -0.13%  ' },
    @{ Row = 39; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.6495'; E = '  +1.56%  ' },
    @{ Row = 40; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.993'; E = '  +2.55%  ' },
    @{ Row = 41; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.226'; E = '  -1.13%  ' },
    @{ Row = 42; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '11.31'; E = '  +1.11%  ' },
    @{ Row = 43; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.6142'; E = '  +1.03%  ' },
    @{ Row = 44; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '13.06'; E = '  +0.19%  ' },
    @{ Row = 45; B = 'WEMIXTOKEN'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '1.285'; E = '  +0.65%  ' },
    @{ Row = 46; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '3.675'; E = '  +0.48%  ' },
    @{ Row = 47; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '2.019'; E = '  +2.28%  ' },
    @{ Row = 48; B = 'EOS'; C = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; D = '1.231'; E = '  +2.46%  ' },
    @{ Row = 49; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '120.81'; E = '  +0.14%  ' },
    @{ Row = 50; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.06911'; E = '  +0.97%  ' },
    @{ Row = 51; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '78.13'; E = '  -1.19%  ' }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        # Prefix with an apostrophe so Excel stores the value as text rather
        # than re-interpreting it as a number/date, then strip the resulting
        # "Text" number format so the cell keeps its original (default) style.
        $ws.Cells.Item($r, 4).Value = "'" + $u.D
        $ws.Cells.Item($r, 4).ClearFormats()
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
